$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text (many values look numeric, e.g. "590.22" or
# "68.310.44", and Excel would silently coerce them to numbers/dates on
# assignment). We set the number format to Text first, assign the value, then
# restore the default "Normal" style so no stray style index is left behind.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.310.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.697.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.688.93"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.22%  "

# Row 8
$ws.Range("E8").Value = "  -6.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.721"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.76%  "

# Row 11
$ws.Range("E11").Value = "  -5.72%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.81%  "

# Row 13
$ws.Range("E13").Value = "  -9.12%  "

# Row 14
$ws.Range("E14").Value = "  -6.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.277.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.698.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.76%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.90%  "

# Row 18
$ws.Range("E18").Value = "  -2.40%  "

# Row 19
$ws.Range("E19").Value = "  -6.09%  "

# Row 20
$ws.Range("E20").Value = "  -6.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.069.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "409.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.80%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.33%  "

# Row 26
$ws.Range("E26").Value = "  -6.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.22%  "

# Row 29
$ws.Range("E29").Value = "  +2.43%  "

# Row 30
$ws.Range("E30").Value = "  -6.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.74%  "

# Row 33
$ws.Range("E33").Value = "  -6.84%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.74%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "590.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.79%  "

# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0437"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.92%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.10%  "

# Row 47
$ws.Range("E47").Value = "  -7.28%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.805.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "

# Row 49
$ws.Range("E49").Value = "  -6.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.21%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.10%  "
